# Applies the scheduled-runner market-price refresh to the leve profit sheets.
# For each affected row, currentAveragePrice[/NQ/HQ], LevePrice[NQ/HQ] and
# LeveProfit[NQ/HQ] (columns H-N) are refreshed to the latest computed values.
# Some rows gain/lose a profit cell entirely when an item flips between having
# a computed profit and being blank (ClearContents mirrors a truly empty cell).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H42").Value = 356.2
$ws.Range("I42").Value = 430.375
$ws.Range("J42").Value = 59.5
$ws.Range("K42").Value = 1291.125
$ws.Range("L42").Value = 178.5
$ws.Range("M42").Value = -1061.125
$ws.Range("N42").Value = -638.5

$ws.Range("H87").Value = 104999.5
$ws.Range("J87").Value = 104999.5
$ws.Range("L87").Value = 104999.5
$ws.Range("N87").Value = -107495.5

$ws.Range("H90").Value = 104999.5
$ws.Range("J90").Value = 104999.5
$ws.Range("L90").Value = 314998.5
$ws.Range("N90").Value = -327478.5

$ws.Range("H106").Value = 4176.3335
$ws.Range("I106").Value = 4198.375
$ws.Range("K106").Value = 4198.375
$ws.Range("M106").Value = -3567.375

$ws.Range("H135").Value = 19187.857
$ws.Range("I135").Value = 601.2273
$ws.Range("K135").Value = 5411.045700000001
$ws.Range("M135").Value = -2876.045700000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14453.524
$ws.Range("I32").Value = 14793.557
$ws.Range("K32").Value = 14793.557
$ws.Range("M32").Value = -14506.557

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 28585
$ws.Range("J110").Value = 28585
$ws.Range("L110").Value = 28585
$ws.Range("N110").Value = -36765

$ws.Range("H134").Value = 2363.375
$ws.Range("I134").Value = 1947.4193
$ws.Range("K134").Value = 5842.257900000001
$ws.Range("M134").Value = -3307.257900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1095
$ws.Range("I22").Value = 674
$ws.Range("J22").Value = 1600.2
$ws.Range("K22").Value = 674
$ws.Range("L22").Value = 1600.2
$ws.Range("M22").Value = -324
$ws.Range("N22").Value = -2300.2

$ws.Range("H37").Value = 40000
$ws.Range("J37").Value = 40000
$ws.Range("L37").Value = 40000
$ws.Range("N37").Value = -40214

$ws.Range("H62").Value = 7718.08
$ws.Range("J62").Value = 7961
$ws.Range("L62").Value = 7961
$ws.Range("N62").Value = -9209

$ws.Range("H65").Value = 7718.08
$ws.Range("J65").Value = 7961
$ws.Range("L65").Value = 39805
$ws.Range("N65").Value = -46045

$ws.Range("H99").Value = 9707.467000000001
$ws.Range("I99").Value = 4868.4
$ws.Range("K99").Value = 4868.4
$ws.Range("M99").Value = -3370.4

$ws.Range("H105").Value = 1615.6666
$ws.Range("I105").Value = 938.8
$ws.Range("K105").Value = 938.8
$ws.Range("M105").Value = 808.2

$ws.Range("H122").Value = 2332.1143
$ws.Range("I122").Value = 2326
$ws.Range("J122").Value = 2361.6667
$ws.Range("K122").Value = 6978
$ws.Range("L122").Value = 7085.000100000001
$ws.Range("M122").Value = -4528
$ws.Range("N122").Value = -11985.0001

$ws.Range("H126").Value = 9707.467000000001
$ws.Range("I126").Value = 4868.4
$ws.Range("K126").Value = 14605.2
$ws.Range("M126").Value = -12135.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 500
$ws.Range("K14").Value = 1500
$ws.Range("M14").Value = -1327

$ws.Range("H23").Value = 550.3
$ws.Range("J23").Value = 634.3570999999999
$ws.Range("L23").Value = 1903.0713
$ws.Range("N23").Value = -2373.0713

$ws.Range("H56").Value = 111111110
$ws.Range("I56").Value = 111111110
$ws.Range("K56").Value = 111111110
$ws.Range("M56").Value = -111110580

$ws.Range("H117").Value = 1215.5714
$ws.Range("I117").Value = 1127.25
$ws.Range("J117").Value = 1333.3334
$ws.Range("K117").Value = 3381.75
$ws.Range("L117").Value = 4000.0002
$ws.Range("M117").Value = 60.25
$ws.Range("N117").Value = -10884.0002

$ws.Range("H121").Value = 64517.883
$ws.Range("I121").Value = 119753.11
$ws.Range("J121").Value = 2378.25
$ws.Range("K121").Value = 359259.33
$ws.Range("L121").Value = 7134.75
$ws.Range("M121").Value = -357949.33
$ws.Range("N121").Value = -9754.75

$ws.Range("H127").Value = 2632.25
$ws.Range("J127").Value = 3166.3333
$ws.Range("L127").Value = 9498.999899999999
$ws.Range("N127").Value = -19418.9999

$ws.Range("H131").Value = 204178.23
$ws.Range("I131").Value = 608459.1
$ws.Range("J131").Value = 2037.7858
$ws.Range("K131").Value = 1825377.3
$ws.Range("L131").Value = 6113.357400000001
$ws.Range("M131").Value = -1820337.3
$ws.Range("N131").Value = -16193.3574

$ws.Range("H134").Value = 2154.8572
$ws.Range("I134").Value = 2130.8333
$ws.Range("J134").Value = 2299
$ws.Range("K134").Value = 6392.499899999999
$ws.Range("L134").Value = 6897
$ws.Range("M134").Value = -1322.499899999999
$ws.Range("N134").Value = -17037

$ws.Range("H139").Value = 6875.75
$ws.Range("I139").Value = 7808.7
$ws.Range("K139").Value = 23426.1
$ws.Range("M139").Value = -18286.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2143.0908
$ws.Range("I126").Value = 1376.8
$ws.Range("J126").Value = 3785.1428
$ws.Range("K126").Value = 4130.4
$ws.Range("L126").Value = 11355.4284
$ws.Range("M126").Value = -1660.4
$ws.Range("N126").Value = -16295.4284

$ws.Range("H132").Value = 2665.0977
$ws.Range("I132").Value = 2617.8918
$ws.Range("J132").Value = 3101.75
$ws.Range("K132").Value = 7853.6754
$ws.Range("L132").Value = 9305.25
$ws.Range("M132").Value = -5323.6754
$ws.Range("N132").Value = -14365.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 19700
$ws.Range("J20").Value = 19700
$ws.Range("L20").Value = 19700
$ws.Range("N20").Value = -20152

$ws.Range("H40").Value = 2363.375
$ws.Range("I40").Value = 1939.3334
$ws.Range("K40").Value = 1939.3334
$ws.Range("M40").Value = -1803.3334

$ws.Range("H46").Value = 4234.3335
$ws.Range("I46").Value = 724.6667
$ws.Range("J46").Value = 7744
$ws.Range("K46").Value = 724.6667
$ws.Range("L46").Value = 7744
$ws.Range("M46").Value = -536.6667
$ws.Range("N46").Value = -8120

$ws.Range("H61").Value = 1498.5
$ws.Range("I61").Value = 998.3
$ws.Range("K61").Value = 998.3
$ws.Range("M61").Value = -796.3

$ws.Range("H113").Value = 1498.5
$ws.Range("I113").Value = 998.3
$ws.Range("K113").Value = 998.3
$ws.Range("M113").Value = 1171.7

$ws.Range("H132").Value = 3480.4285
$ws.Range("I132").Value = 3480.4285
$ws.Range("K132").Value = 10441.2855
$ws.Range("M132").Value = -7911.2855

$ws.Range("H136").Value = 3833.3333
$ws.Range("I136").Value = 3750
$ws.Range("K136").Value = 11250
$ws.Range("M136").Value = -8700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 16499.666
$ws.Range("I28").Value = 9500
$ws.Range("J28").Value = 19999.5
$ws.Range("K28").Value = 9500
$ws.Range("L28").Value = 19999.5
$ws.Range("M28").Value = -9152
$ws.Range("N28").Value = -20695.5

$ws.Range("H33").Value = 21566.666
$ws.Range("J33").Value = 29850
$ws.Range("L33").Value = 29850
$ws.Range("N33").Value = -30350

$ws.Range("H36").Value = 21566.666
$ws.Range("J36").Value = 29850
$ws.Range("L36").Value = 29850
$ws.Range("N36").Value = -30350

$ws.Range("H122").Value = 39252.37
$ws.Range("I122").Value = 44781.133
$ws.Range("J122").Value = 6079.8
$ws.Range("K122").Value = 134343.399
$ws.Range("L122").Value = 18239.4
$ws.Range("M122").Value = -131893.399
$ws.Range("N122").Value = -23139.4

$ws.Range("H126").Value = 4503.4
$ws.Range("J126").Value = 5095.636
$ws.Range("L126").Value = 15286.908
$ws.Range("N126").Value = -20226.908

$ws.Range("H132").Value = 40770.168
$ws.Range("I132").Value = 48190.934
$ws.Range("K132").Value = 144572.802
$ws.Range("M132").Value = -142042.802
